$wb = $excel.ActiveWorkbook

$about  = $wb.Worksheets.Item("About")
$ws2    = $wb.Worksheets.Item("CDCF-PMpPDOU")
$ws3    = $wb.Worksheets.Item("CDCF-FTMpFDOU")

# ---------------------------------------------------------------------------
# "About" sheet: update text to reflect the India model (km instead of miles)
# and add a "Google unit converter" hyperlink as the Source.
#
# Shift the "Notes" block (and everything below it) down by one row to make
# room for the two new hyperlink rows (5 & 6).
# ---------------------------------------------------------------------------
$about.Rows("6").Insert()

$about.Range("A1").Value = "CDCF Passenger Miles per Passenger Distance Output Unit"
$about.Range("A2").Value = "CDCF Freight Ton Miles per Freight Distance Output Unit"

$about.Range("A4").Value = "Source:"
$about.Range("B4").Value = "Google unit converter"

# B5 gets the hyperlink; B6 is left as an empty, hyperlink-styled cell
# (mirrors what Excel leaves behind when the link is followed by a blank,
# similarly-styled row underneath it).
$about.Hyperlinks.Add($about.Range("B5"), "https://www.google.com/search?q=kilometers+per+mile", "", "", "https://www.google.com/search?q=kilometers+per+mile") | Out-Null
$about.Range("B6").Style = "Hyperlink"

$about.Range("A7").Value = "Notes"
$about.Range("A8").Value = "This variable converts the cargo distance units used by the model"
$about.Range("A9").Value = "internally (passenger*miles) and (freight ton*miles) to the"
$about.Range("A10").Value = "desired output unit."

$about.Range("A12").Value = "For the India model, the desired output units are:"
$about.Range("A13").Value = "trillion passenger-kilometers"
$about.Range("A14").Value = "trillion freight ton-kilometers"

# ---------------------------------------------------------------------------
# "CDCF-PMpPDOU" sheet: passenger-miles -> passenger-kilometers conversion
# ---------------------------------------------------------------------------
$ws2.Range("A2").Value = "passenger-kilometers"
$ws2.Range("B2").Formula = "=1.60934*10^12"

# ---------------------------------------------------------------------------
# "CDCF-FTMpFDOU" sheet: freight ton-miles -> freight ton-kilometers conversion
# ---------------------------------------------------------------------------
$ws3.Range("A2").Value = "freight ton-kilometers"
$ws3.Range("B2").Formula = "=1.60934*10^12"
